# In progress: build combinator generator that issues all possible combinations
# of "0"s and "1"s for a given tuple of clues (either for a row or column).
#
# This session appends a freshly generated sample grid (a column-index header
# row plus five generated 0/1 combinations) below the existing nonogram
# solver sections, and nudges the sheet/window view to where the new data
# lives.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column widths: widen the clue columns a touch and give the newly used
# columns (I through AL) the same narrow "grid cell" width.
# ---------------------------------------------------------------------------
$ws.Columns("B:F").ColumnWidth = 2.3
$ws.Columns("I:AA").ColumnWidth = 2.3
$ws.Columns("AB:AL").ColumnWidth = 2.6

# ---------------------------------------------------------------------------
# Row 28: a header strip enumerating the column index (1, 2, 3, ...) above
# the generated combinations - the generator's "ruler".
# ---------------------------------------------------------------------------
$headerCols = @("I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z", `
                "AA","AB","AC","AD","AE","AF","AG","AH","AI")
$n = 0
foreach ($c in $headerCols) {
    $n = $n + 1
    $ws.Range($c + "28").Value = $n
}

# ---------------------------------------------------------------------------
# Rows 30-34: five generated 0/1 combinations (only the "1" cells are
# written - blanks stand in for "0", matching the solver's sparse style
# used throughout the rest of the sheet).
# ---------------------------------------------------------------------------
$generatedRows = @{
    30 = @("I","M","Q","U","V","W","X","Z","AD","AF","AG","AH","AI")
    31 = @("I","J","L","M","P","R","U","Z","AD","AF")
    32 = @("I","K","M","O","S","U","V","W","AA","AC","AF","AG","AH")
    33 = @("I","M","O","P","Q","R","S","U","AA","AC","AF")
    34 = @("I","M","O","S","U","V","W","X","AB","AF","AG","AH","AI")
}

foreach ($r in 30..34) {
    foreach ($c in $generatedRows[$r]) {
        $ws.Range($c + $r).Value = 1
    }
}

# ---------------------------------------------------------------------------
# View state: scroll down to the newly generated block and land the
# selection where the generator cursor ended up.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("AM32").Select()
